$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 87979
$ws.Range("A3").Value = 87980
$ws.Range("A4").Value = 87981
$ws.Range("A5").Value = 87982
$ws.Range("A6").Value = 87983
$ws.Range("A7").Value = 87984
$ws.Range("A8").Value = 87985

$ws.Range("C10").Select()
